# Word COM-interop script: "add finding mutual friends problem"
#
# This recreates a set of benchmark-run notes for the mutualFriends
# problem: the original three timing lines get their split
# runs/proofErr-marker artifacts cleaned up into single plain runs, and
# two new result blocks (each a label paragraph followed by a raw JSON
# metrics paragraph) are appended around the existing blank paragraph.

$d = $word.ActiveDocument

# --- Step 1: normalize the first three paragraphs -------------------------
# Each of these paragraphs originally held its text split across three
# runs with <w:proofErr w:type="gramStart"/>/<w:proofErr w:type="gramEnd"/>
# markers sandwiched in between (left over from the grammar checker).
# Running Find/Replace against the paragraph's plain text collapses it
# back down to a single run and drops the proofErr markers.
$d.Content.Find.Execute("100 * 100 : 16.989s", $false, $false, $false, $false, $false, $true, 1, $false, "100 * 100 : 16.989s", 2) | Out-Null
$d.Content.Find.Execute("500 * 500 : 20.94s", $false, $false, $false, $false, $false, $true, 1, $false, "500 * 500 : 20.94s", 2) | Out-Null
$d.Content.Find.Execute("700 *700 : 1.19.688 m", $false, $false, $false, $false, $false, $true, 1, $false, "700 *700 : 1.19.688 m", 2) | Out-Null

# --- Step 2: insert the "3 tabs" result block ------------------------------
# Goes right after the bookmark paragraph (paragraph 4) and before the
# pre-existing blank paragraph (paragraph 5), without touching that blank
# paragraph itself.
$afterBookmark = $d.Paragraphs(4).Range
$afterBookmark.InsertParagraphAfter()
$afterBookmark.InsertParagraphAfter()

$d.Paragraphs(5).Range.Text = "3 tabs with 20 userIDs"
$d.Paragraphs(6).Range.Text = '{"description":"mapping process is started at 2021-11-25 13:51:09 - mapping is finished and reducing has started. date: 2021-11-25 13:51:34 - owner job process has completed at: 2021-11-25 13:51:44","total_ownerJob_duration":36.207276,"transformed_data_size":15664,"request_count":26,"response_count":25,"total_server_process":18.20694,"server_process_duration_time_detail":"0.378621,0.314024,1.163829,1.311038,1.363816,0.225657,0.300541,0.26054,0.364043,0.261538,0.287815,0.21763,0.283109,0.27199,0.279784,0.767362,0.292852,0.300857,0.235264,0.239661,0.357822,0.31887,0.240492,0.205779,10.817494","metadata_request_size":35048,"metadata_response_size":30000,"total_ocuupied_bandwidth":80744}'

# --- Step 3: insert the "1 tabs" result block ------------------------------
# Goes at the very end of the document, after the pre-existing blank
# paragraph (now paragraph 7), again without touching that paragraph.
$d.Range($d.Content.End, $d.Content.End).InsertParagraphAfter()
$d.Range($d.Content.End, $d.Content.End).InsertParagraphAfter()

$d.Paragraphs(8).Range.Text = "1 tabs with 20 userIDs:"
$d.Paragraphs(9).Range.Text = '{"description":"mapping process is started at 2021-11-25 13:59:15 - mapping is finished and reducing has started. date: 2021-11-25 13:59:51 - owner job process has completed at: 2021-11-25 14:00:02","total_ownerJob_duration":47.141443,"transformed_data_size":14555,"request_count":21,"response_count":21,"total_server_process":17.275978000000002,"server_process_duration_time_detail":"0.337309,1.249752,0.266253,0.287478,0.233893,0.262301,0.289867,0.273468,0.231479,0.24069,0.287933,0.276308,0.291437,0.261106,0.237839,0.329673,0.232826,0.364252,0.267661,0.230197,10.824256","metadata_request_size":28308,"metadata_response_size":25200,"total_ocuupied_bandwidth":68095}'
